# Domain experts (primary reviewer) shifted so Bo can review double pend,
# some associated shifts in secondary reviewers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Primary "Expert" (column J) swap -----------------------------------
# Row 7 (Sasha Soraine's project): expert Bo Cao -> Peter Michalski
# Row 9 (Zhi Zhang's "Double pendulum" project): expert Peter Michalski -> Bo Cao
# (this frees Bo Cao up to review the double pendulum project)
$ws.Range("J7").Value = "Peter Michalski"
$ws.Range("J9").Value = "Bo Cao"

# --- Secondary reviewer chain (columns L/M) updates ----------------------
$ws.Range("L6").Formula = "=K3"
$ws.Range("L9").Formula = "=K7"
$ws.Range("M5").Formula = "=L9"
$ws.Range("M8").Formula = "=L6"

# --- Selection / active cell ---------------------------------------------
$ws.Range("M9").Select()
